# Regenerate the "K" column (G) values for each outing row.
# The sheet's column G header is "K" (strikeouts), which previously held a
# different derived quantity ("Strike#"). This recalculates/overwrites the
# per-row K values with the corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    3  = 1
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 2
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 0
    18 = 2
    19 = 4
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 2
    25 = 0
    26 = 1
    27 = 2
    28 = 0
    29 = 2
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 3
    38 = 1
    39 = 2
    40 = 0
    41 = 1
    42 = 2
    43 = 0
    44 = 2
    45 = 3
    46 = 0
    47 = 2
    48 = 1
    49 = 2
    50 = 1
    51 = 0
    52 = 2
    53 = 1
    54 = 2
    55 = 2
    56 = 2
    57 = 0
    58 = 3
    59 = 1
    60 = 2
    61 = 0
    62 = 2
    64 = 1
    65 = 2
    66 = 1
    67 = 2
    68 = 5
    69 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
